$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.118.61'
$ws.Range("E2").Value = '  +2.91%  '

# Row 3
$ws.Range("D3").Value = '1.578.16'
$ws.Range("E3").Value = '  +1.78%  '

# Row 4
$ws.Range("E4").Value = '  -0.55%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.86%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.54%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '25.86'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.73%  '

# Row 9
$ws.Range("E9").Value = '  +2.57%  '

# Row 10
$ws.Range("E10").Value = '  +1.68%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.18%  '

# Row 12
$ws.Range("D12").Value = '1.803.98'
$ws.Range("E12").Value = '  +1.86%  '

# Row 13
$ws.Range("D13").Value = '1.586.74'
$ws.Range("E13").Value = '  +2.35%  '

# Row 14
$ws.Range("D14").Value = '29.116.82'
$ws.Range("E14").Value = '  +2.94%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.521'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.43%  '

# Row 16
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.07%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.93%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.57%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0691'
$ws.Range("E20").Value = '  +2.71%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.98'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.88%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.12%  '

# Row 26
$ws.Range("E26").Value = '  +4.38%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.45%  '

# Row 28
$ws.Range("E28").Value = '  +1.33%  '

# Row 29
$ws.Range("E29").Value = '  -0.55%  '

# Row 30
$ws.Range("E30").Value = '  -0.39%  '

# Row 31
$ws.Range("E31").Value = '  +0.10%  '

# Row 32
$ws.Range("E32").Value = '  +1.34%  '

# Row 33
$ws.Range("D33").Value = '1.419.41'
$ws.Range("E33").Value = '  +2.71%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.05'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.30%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.04'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.18%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.60%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.23%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.28'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.22%  '

# Row 39
$ws.Range("E39").Value = '  +1.52%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.525'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.52%  '

# Row 41
$ws.Range("E41").Value = '  +2.11%  '

# Row 42
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.54%  '

# Row 43
$ws.Range("B43").Value = 'BitcoinSV'
$ws.Range("C43").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '52.40'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +23.76%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.788'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.66%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0473'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.59%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.85%  '

# Row 48
$ws.Range("D48").Value = '1.716.23'
$ws.Range("E48").Value = '  +1.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.848'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '85.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0513'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
